$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 3 holds the 7527ec17... entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-20 02:33:45"
$wsZhCn.Range("H3").Value = "2016-03-20 02:34:05"

# de-de sheet: row 3 holds the 7527ec17... entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-20 02:33:48"
$wsDeDe.Range("H3").Value = "2016-03-20 02:34:10"
